$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing header/body cells (column C) onto the new
# column D cells first, so the cellXfs styles 1 (header) and 2 (body) get
# reused rather than new styles being minted.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("D2:D16").PasteSpecial(-4122)

# Now populate the values. "up" is entered first (so it lands earlier in
# the shared-string table), then "down", then the header "direction" last.
$ws.Range("D2").Value = "up"
$ws.Range("D3").Value = "down"
$ws.Range("D4").Value = "down"
$ws.Range("D5").Value = "down"
$ws.Range("D6").Value = "down"
$ws.Range("D7").Value = "down"
$ws.Range("D8").Value = "down"
$ws.Range("D9").Value = "down"
$ws.Range("D10").Value = "down"
$ws.Range("D11").Value = "down"
$ws.Range("D12").Value = "up"
$ws.Range("D13").Value = "up"
$ws.Range("D14").Value = "down"
$ws.Range("D15").Value = "down"
$ws.Range("D16").Value = "down"
$ws.Range("D1").Value = "direction"

# Match the final cursor position recorded in the workbook.
$ws.Range("D2").Select()
